# Refactor/bugfix commit: rename "PtLabresult" -> "PtLabResult" and switch
# the active/selected tab from the first sheet (PtAssessment) to the
# third sheet (PtLabResult).

$wb = $excel.ActiveWorkbook

# Third sheet: "PtLabresult" -> fix the casing to "PtLabResult".
$wsLabResult = $wb.Worksheets.Item(3)
$wsLabResult.Name = "PtLabResult"

# Make it the active/selected sheet (moves tabSelected from sheet 1 to
# sheet 3, and updates the workbook's bookViews/activeTab accordingly).
$wsLabResult.Activate()
